# Update LR-pairs TPM-derived values (Sema3d-Nrp1) with new TPM-based computations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.5890554026108095
$ws.Range("J2").Value = 0.5890554026108095
$ws.Range("M2").Value = 86.89540866666668
$ws.Range("N2").Value = 260.686226
$ws.Range("O2").Value = 0.319779657009892
$ws.Range("P2").Value = 0.3197796570098919
$ws.Range("Q2").Value = 91.96022342134823
$ws.Range("R2").Value = 827.642010792134
$ws.Range("S2").Value = 0.1883679346067085
$ws.Range("T2").Value = 0.1883679346067084

# Row 3
$ws.Range("I3").Value = 0.5890554026108095
$ws.Range("J3").Value = 0.5890554026108095
$ws.Range("O3").Value = 0.1999969065479545
$ws.Range("P3").Value = 0.1999969065479545
$ws.Range("S3").Value = 0.1178092583075218
$ws.Range("T3").Value = 0.1178092583075218

# Row 4
$ws.Range("I4").Value = 0.5890554026108095
$ws.Range("J4").Value = 0.5890554026108095
$ws.Range("M4").Value = 60.92601633333334
$ws.Range("N4").Value = 182.778049
$ws.Range("O4").Value = 0.224210932487692
$ws.Range("P4").Value = 0.224210932487692
$ws.Range("Q4").Value = 64.4771704300101
$ws.Range("R4").Value = 580.294533870091
$ws.Range("S4").Value = 0.1320726611062824
$ws.Range("T4").Value = 0.1320726611062824

# Row 5
$ws.Range("I5").Value = 0.5890554026108095
$ws.Range("J5").Value = 0.5890554026108095
$ws.Range("M5").Value = 7.809668333333332
$ws.Range("N5").Value = 23.429005
$ws.Range("O5").Value = 0.02873998867505581
$ws.Range("P5").Value = 0.02873998867505581
$ws.Range("Q5").Value = 8.264865265032775
$ws.Range("R5").Value = 74.38378738529498
$ws.Range("S5").Value = 0.01692944560001511
$ws.Range("T5").Value = 0.01692944560001511

# Row 6
$ws.Range("I6").Value = 0.5890554026108095
$ws.Range("J6").Value = 0.5890554026108095
$ws.Range("M6").Value = 61.75795633333333
$ws.Range("N6").Value = 185.273869
$ws.Range("O6").Value = 0.2272725152794058
$ws.Range("P6").Value = 0.2272725152794058
$ws.Range("Q6").Value = 65.35760116216343
$ws.Range("R6").Value = 588.2184104594709
$ws.Range("S6").Value = 0.1338761029902817
$ws.Range("T6").Value = 0.1338761029902817

# Row 7
$ws.Range("G7").Value = 0.7382956666666667
$ws.Range("H7").Value = 2.214887
$ws.Range("I7").Value = 0.4109445973891905
$ws.Range("J7").Value = 0.4109445973891905
$ws.Range("M7").Value = 86.89540866666668
$ws.Range("N7").Value = 260.686226
$ws.Range("O7").Value = 0.319779657009892
$ws.Range("P7").Value = 0.3197796570098919
$ws.Range("Q7").Value = 64.15450367182912
$ws.Range("R7").Value = 577.3905330464621
$ws.Range("S7").Value = 0.1314117224031835
$ws.Range("T7").Value = 0.1314117224031835

# Row 8
$ws.Range("G8").Value = 0.7382956666666667
$ws.Range("H8").Value = 2.214887
$ws.Range("I8").Value = 0.4109445973891905
$ws.Range("J8").Value = 0.4109445973891905
$ws.Range("O8").Value = 0.1999969065479545
$ws.Range("P8").Value = 0.1999969065479545
$ws.Range("Q8").Value = 40.12357257324945
$ws.Range("R8").Value = 361.112153159245
$ws.Range("S8").Value = 0.08218764824043275
$ws.Range("T8").Value = 0.08218764824043273

# Row 9
$ws.Range("G9").Value = 0.7382956666666667
$ws.Range("H9").Value = 2.214887
$ws.Range("I9").Value = 0.4109445973891905
$ws.Range("J9").Value = 0.4109445973891905
$ws.Range("M9").Value = 60.92601633333334
$ws.Range("N9").Value = 182.778049
$ws.Range("O9").Value = 0.224210932487692
$ws.Range("P9").Value = 0.224210932487692
$ws.Range("Q9").Value = 44.98141384616256
$ws.Range("R9").Value = 404.832724615463
$ws.Range("S9").Value = 0.09213827138140956
$ws.Range("T9").Value = 0.09213827138140956

# Row 10
$ws.Range("G10").Value = 0.7382956666666667
$ws.Range("H10").Value = 2.214887
$ws.Range("I10").Value = 0.4109445973891905
$ws.Range("J10").Value = 0.4109445973891905
$ws.Range("M10").Value = 7.809668333333332
$ws.Range("N10").Value = 23.429005
$ws.Range("O10").Value = 0.02873998867505581
$ws.Range("P10").Value = 0.02873998867505581
$ws.Range("Q10").Value = 5.765844288603888
$ws.Range("R10").Value = 51.892598597435
$ws.Range("S10").Value = 0.01181054307504071
$ws.Range("T10").Value = 0.01181054307504071

# Row 11
$ws.Range("G11").Value = 0.7382956666666667
$ws.Range("H11").Value = 2.214887
$ws.Range("I11").Value = 0.4109445973891905
$ws.Range("J11").Value = 0.4109445973891905
$ws.Range("M11").Value = 61.75795633333333
$ws.Range("N11").Value = 185.273869
$ws.Range("O11").Value = 0.2272725152794058
$ws.Range("P11").Value = 0.2272725152794058
$ws.Range("Q11").Value = 45.59563154308923
$ws.Range("R11").Value = 410.360683887803
$ws.Range("S11").Value = 0.09339641228912408
$ws.Range("T11").Value = 0.09339641228912407
